$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 6.112374333333332
$ws.Range("H2").Value = 18.337123
$ws.Range("I2").Value = 0.02125075796073207
$ws.Range("J2").Value = 0.02125075796073207
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 14.08030833333333
$ws.Range("N2").Value = 42.240925
$ws.Range("O2").Value = 0.1556474304006564
$ws.Range("P2").Value = 0.1556474304006564
$ws.Range("Q2").Value = 86.06411526208609
$ws.Range("R2").Value = 774.5770373587749
$ws.Range("S2").Value = 0.00330762587065424
$ws.Range("T2").Value = 0.003307625870654241

# Row 3
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 6.112374333333332
$ws.Range("H3").Value = 18.337123
$ws.Range("I3").Value = 0.02125075796073207
$ws.Range("J3").Value = 0.02125075796073207
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 43.15322233333333
$ws.Range("N3").Value = 129.459667
$ws.Range("O3").Value = 0.4770270657916382
$ws.Range("P3").Value = 0.4770270657916382
$ws.Range("Q3").Value = 263.7686485908934
$ws.Range("R3").Value = 2373.917837318041
$ws.Range("S3").Value = 0.01013718671585632
$ws.Range("T3").Value = 0.01013718671585632

# Row 4
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 6.112374333333332
$ws.Range("H4").Value = 18.337123
$ws.Range("I4").Value = 0.02125075796073207
$ws.Range("J4").Value = 0.02125075796073207
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 21.909414
$ws.Range("N4").Value = 65.72824199999999
$ws.Range("O4").Value = 0.2421924229181179
$ws.Range("P4").Value = 0.242192422918118
$ws.Range("Q4").Value = 133.918539791974
$ws.Range("R4").Value = 1205.266858127766
$ws.Range("S4").Value = 0.005146772559356184
$ws.Range("T4").Value = 0.005146772559356184

# Row 5
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 6.112374333333332
$ws.Range("H5").Value = 18.337123
$ws.Range("I5").Value = 0.02125075796073207
$ws.Range("J5").Value = 0.02125075796073207
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 11.31989366666667
$ws.Range("N5").Value = 33.959681
$ws.Range("O5").Value = 0.1251330808895874
$ws.Range("P5").Value = 0.1251330808895874
$ws.Range("Q5").Value = 69.19142750419589
$ws.Range("R5").Value = 622.722847537763
$ws.Range("S5").Value = 0.002659172814865329
$ws.Range("T5").Value = 0.00265917281486533

# Row 6
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 263.0904693333333
$ws.Range("H6").Value = 789.271408
$ws.Range("I6").Value = 0.9146808720612395
$ws.Range("J6").Value = 0.9146808720612395
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 14.08030833333333
$ws.Range("N6").Value = 42.240925
$ws.Range("O6").Value = 0.1556474304006564
$ws.Range("P6").Value = 0.1556474304006564
$ws.Range("Q6").Value = 3704.394927774711
$ws.Range("R6").Value = 33339.55434997239
$ws.Range("S6").Value = 0.1423677273729635
$ws.Range("T6").Value = 0.1423677273729635

# Row 7
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 263.0904693333333
$ws.Range("H7").Value = 789.271408
$ws.Range("I7").Value = 0.9146808720612395
$ws.Range("J7").Value = 0.9146808720612395
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 43.15322233333333
$ws.Range("N7").Value = 129.459667
$ws.Range("O7").Value = 0.4770270657916382
$ws.Range("P7").Value = 0.4770270657916382
$ws.Range("Q7").Value = 11353.20151692235
$ws.Range("R7").Value = 102178.8136523011
$ws.Range("S7").Value = 0.4363275325351099
$ws.Range("T7").Value = 0.4363275325351099

# Row 8
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 263.0904693333333
$ws.Range("H8").Value = 789.271408
$ws.Range("I8").Value = 0.9146808720612395
$ws.Range("J8").Value = 0.9146808720612395
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 21.909414
$ws.Range("N8").Value = 65.72824199999999
$ws.Range("O8").Value = 0.2421924229181179
$ws.Range("P8").Value = 0.242192422918118
$ws.Range("Q8").Value = 5764.158012078303
$ws.Range("R8").Value = 51877.42210870473
$ws.Range("S8").Value = 0.2215287766013686
$ws.Range("T8").Value = 0.2215287766013687

# Row 9
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 263.0904693333333
$ws.Range("H9").Value = 789.271408
$ws.Range("I9").Value = 0.9146808720612395
$ws.Range("J9").Value = 0.9146808720612395
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 11.31989366666667
$ws.Range("N9").Value = 33.959681
$ws.Range("O9").Value = 0.1251330808895874
$ws.Range("P9").Value = 0.1251330808895874
$ws.Range("Q9").Value = 2978.156137566761
$ws.Range("R9").Value = 26803.40523810085
$ws.Range("S9").Value = 0.1144568355517974
$ws.Range("T9").Value = 0.1144568355517974

# Row 10
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 8.317515
$ws.Range("H10").Value = 24.952545
$ws.Range("I10").Value = 0.02891732221566466
$ws.Range("J10").Value = 0.02891732221566466
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 14.08030833333333
$ws.Range("N10").Value = 42.240925
$ws.Range("O10").Value = 0.1556474304006564
$ws.Range("P10").Value = 0.1556474304006564
$ws.Range("Q10").Value = 117.113175767125
$ws.Range("R10").Value = 1054.018581904125
$ws.Range("S10").Value = 0.004500906896936019
$ws.Range("T10").Value = 0.00450090689693602

# Row 11
$ws.Range("E11").Value = 3
$ws.Range("F11").Value = 1
$ws.Range("G11").Value = 8.317515
$ws.Range("H11").Value = 24.952545
$ws.Range("I11").Value = 0.02891732221566466
$ws.Range("J11").Value = 0.02891732221566466
$ws.Range("K11").Value = 3
$ws.Range("L11").Value = 1
$ws.Range("M11").Value = 43.15322233333333
$ws.Range("N11").Value = 129.459667
$ws.Range("O11").Value = 0.4770270657916382
$ws.Range("P11").Value = 0.4770270657916382
$ws.Range("Q11").Value = 358.927574055835
$ws.Range("R11").Value = 3230.348166502515
$ws.Range("S11").Value = 0.01379434536708986
$ws.Range("T11").Value = 0.01379434536708986

# Row 12
$ws.Range("E12").Value = 3
$ws.Range("F12").Value = 1
$ws.Range("G12").Value = 8.317515
$ws.Range("H12").Value = 24.952545
$ws.Range("I12").Value = 0.02891732221566466
$ws.Range("J12").Value = 0.02891732221566466
$ws.Range("K12").Value = 3
$ws.Range("L12").Value = 1
$ws.Range("M12").Value = 21.909414
$ws.Range("N12").Value = 65.72824199999999
$ws.Range("O12").Value = 0.2421924229181179
$ws.Range("P12").Value = 0.242192422918118
$ws.Range("Q12").Value = 182.23187958621
$ws.Range("R12").Value = 1640.08691627589
$ws.Range("S12").Value = 0.007003556331715741
$ws.Range("T12").Value = 0.007003556331715742

# Row 13
$ws.Range("E13").Value = 3
$ws.Range("F13").Value = 1
$ws.Range("G13").Value = 8.317515
$ws.Range("H13").Value = 24.952545
$ws.Range("I13").Value = 0.02891732221566466
$ws.Range("J13").Value = 0.02891732221566466
$ws.Range("K13").Value = 3
$ws.Range("L13").Value = 1
$ws.Range("M13").Value = 11.31989366666667
$ws.Range("N13").Value = 33.959681
$ws.Range("O13").Value = 0.1251330808895874
$ws.Range("P13").Value = 0.1251330808895874
$ws.Range("Q13").Value = 94.15338537090501
$ws.Range("R13").Value = 847.3804683381451
$ws.Range("S13").Value = 0.003618513619923028
$ws.Range("T13").Value = 0.003618513619923029

# Row 14
$ws.Range("E14").Value = 3
$ws.Range("F14").Value = 1
$ws.Range("G14").Value = 10.110527
$ws.Range("H14").Value = 30.331581
$ws.Range("I14").Value = 0.03515104776236379
$ws.Range("J14").Value = 0.03515104776236379
$ws.Range("K14").Value = 3
$ws.Range("L14").Value = 1
$ws.Range("M14").Value = 14.08030833333333
$ws.Range("N14").Value = 42.240925
$ws.Range("O14").Value = 0.1556474304006564
$ws.Range("P14").Value = 0.1556474304006564
$ws.Range("Q14").Value = 142.3593375724917
$ws.Range("R14").Value = 1281.234038152425
$ws.Range("S14").Value = 0.005471170260102668
$ws.Range("T14").Value = 0.005471170260102669

# Row 15
$ws.Range("E15").Value = 3
$ws.Range("F15").Value = 1
$ws.Range("G15").Value = 10.110527
$ws.Range("H15").Value = 30.331581
$ws.Range("I15").Value = 0.03515104776236379
$ws.Range("J15").Value = 0.03515104776236379
$ws.Range("K15").Value = 3
$ws.Range("L15").Value = 1
$ws.Range("M15").Value = 43.15322233333333
$ws.Range("N15").Value = 129.459667
$ws.Range("O15").Value = 0.4770270657916382
$ws.Range("P15").Value = 0.4770270657916382
$ws.Range("Q15").Value = 436.3018195381696
$ws.Range("R15").Value = 3926.716375843527
$ws.Range("S15").Value = 0.01676800117358213
$ws.Range("T15").Value = 0.01676800117358213

# Row 16
$ws.Range("E16").Value = 3
$ws.Range("F16").Value = 1
$ws.Range("G16").Value = 10.110527
$ws.Range("H16").Value = 30.331581
$ws.Range("I16").Value = 0.03515104776236379
$ws.Range("J16").Value = 0.03515104776236379
$ws.Range("K16").Value = 3
$ws.Range("L16").Value = 1
$ws.Range("M16").Value = 21.909414
$ws.Range("N16").Value = 65.72824199999999
$ws.Range("O16").Value = 0.2421924229181179
$ws.Range("P16").Value = 0.242192422918118
$ws.Range("Q16").Value = 221.515721801178
$ws.Range("R16").Value = 1993.641496210602
$ws.Range("S16").Value = 0.008513317425677375
$ws.Range("T16").Value = 0.008513317425677376

# Row 17
$ws.Range("E17").Value = 3
$ws.Range("F17").Value = 1
$ws.Range("G17").Value = 10.110527
$ws.Range("H17").Value = 30.331581
$ws.Range("I17").Value = 0.03515104776236379
$ws.Range("J17").Value = 0.03515104776236379
$ws.Range("K17").Value = 3
$ws.Range("L17").Value = 1
$ws.Range("M17").Value = 11.31989366666667
$ws.Range("N17").Value = 33.959681
$ws.Range("O17").Value = 0.1251330808895874
$ws.Range("P17").Value = 0.1251330808895874
$ws.Range("Q17").Value = 114.4500905539623
$ws.Range("R17").Value = 1030.050814985661
$ws.Range("S17").Value = 0.003618513619923028
$ws.Range("T17").Value = 0.003618513619923029

